$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to write (A1 ref -> new text). These mirror the scraped
# coinranking.com snapshot captured by the scheduled GitHub Actions run.
$changes = [ordered]@{
    "D2" = "58.566.93"
    "E2" = "  -2.70%  "
    "D3" = "2.290.58"
    "E3" = "  -5.39%  "
    "E4" = "  -0.09%  "
    "D5" = "546.82"
    "E5" = "  -1.37%  "
    "D6" = "131.23"
    "E6" = "  -4.42%  "
    "E7" = "  -0.02%  "
    "E8" = "  -2.63%  "
    "D9" = "2.288.43"
    "E9" = "  -5.38%  "
    "E10" = "  -3.14%  "
    "D11" = "5.55"
    "E11" = "  -2.61%  "
    "D12" = "0.149"
    "E12" = "  +0.51%  "
    "E13" = "  -4.92%  "
    "D14" = "23.80"
    "E14" = "  -4.29%  "
    "D15" = "2.696.76"
    "E15" = "  -5.45%  "
    "D16" = "58.520.14"
    "E16" = "  -2.63%  "
    "E17" = "  -3.08%  "
    "D18" = "2.304.09"
    "E18" = "  -5.51%  "
    "D19" = "10.65"
    "E19" = "  -5.45%  "
    "D20" = "4.31"
    "E20" = "  -4.20%  "
    "D21" = "315.46"
    "E21" = "  -3.57%  "
    "D22" = "6.47"
    "E22" = "  -4.04%  "
    "D23" = "1.00"
    "E23" = "  +0.08%  "
    "D24" = "62.84"
    "E24" = "  -3.70%  "
    "E25" = "  -4.06%  "
    "D26" = "0.999"
    "E26" = "  -0.18%  "
    "D27" = "8.14"
    "E27" = "  -7.22%  "
    "D28" = "1.32"
    "E28" = "  -4.77%  "
    "D29" = "1.74"
    "E29" = "  -1.65%  "
    "D30" = "171.27"
    "E30" = "  +0.53%  "
    "D31" = "0.0₃0728"
    "E31" = "  -5.94%  "
    "D32" = "5.81"
    "E32" = "  -4.73%  "
    "E33" = "  -0.23%  "
    "E34" = "  -5.00%  "
    "E35" = "  +0.04%  "
    "D36" = "17.86"
    "E36" = "  -3.66%  "
    "E37" = "  -0.02%  "
    "E38" = "  -5.40%  "
    "D39" = "3.97"
    "E39" = "  -5.66%  "
    "D40" = "37.84"
    "E40" = "  -2.47%  "
    "E41" = "  -5.04%  "
    "D42" = "299.24"
    "E42" = "  -9.04%  "
    "D43" = "140.93"
    "E43" = "  -2.65%  "
    "D44" = "3.46"
    "E44" = "  -5.52%  "
    "D45" = "0.0950"
    "E45" = "  -1.59%  "
    "D46" = "0.0500"
    "E46" = "  -3.01%  "
    "B47" = "InjectiveProtocol"
    "C47" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D47" = "18.60"
    "E47" = "  -7.44%  "
    "B48" = "Mantle"
    "C48" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D48" = "0.555"
    "E48" = "  -3.48%  "
    "E49" = "  -3.43%  "
    "D50" = "16.65"
    "E50" = "  -5.12%  "
    "D51" = "11.01"
    "E51" = "  -0.25%  "
}

# A handful of price cells are pure numeric-looking strings (e.g. "1.00",
# "23.80") that Excel would otherwise coerce to a Number and silently drop
# the trailing zero / significant digits. Force those specific cells to
# keep a Text format so the literal string round-trips unchanged; leave
# every other cell's formatting untouched.
$forceText = @("D14", "D23", "D45", "D46", "D47")

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    if ($forceText -contains $ref) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $changes[$ref]
}
